$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "52.227.65"
$ws.Range("E2").Value = "  +1.02%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.804.87"
$ws.Range("E3").Value = "  +1.83%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "350.62"
$ws.Range("E5").Value = "  +5.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "116.14"
$ws.Range("E6").Value = "  -1.19%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.549"
$ws.Range("E7").Value = "  +2.41%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.593"
$ws.Range("E9").Value = "  +2.88%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.51"
$ws.Range("E10").Value = "  +2.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0862"
$ws.Range("E11").Value = "  +3.67%  "

$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.98"
$ws.Range("E12").Value = "  -1.24%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.131"
$ws.Range("E13").Value = "  +1.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.87"
$ws.Range("E14").Value = "  +3.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.232.78"
$ws.Range("E15").Value = "  +1.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.795.49"
$ws.Range("E16").Value = "  +1.55%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.894"
$ws.Range("E17").Value = "  +0.52%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "52.169.77"
$ws.Range("E18").Value = "  +0.98%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.18"
$ws.Range("E19").Value = "  +5.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.30"
$ws.Range("E20").Value = "  +6.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.43"
$ws.Range("E21").Value = "  -2.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0983"
$ws.Range("E22").Value = "  +1.66%  "

$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.33"
$ws.Range("E23").Value = "  +0.13%  "

$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.10"
$ws.Range("E24").Value = "  -2.89%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.79"
$ws.Range("E25").Value = "  +4.28%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.88"
$ws.Range("E26").Value = "  -0.22%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.26"
$ws.Range("E28").Value = "  -0.79%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.26"
$ws.Range("E29").Value = "  +1.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.141"
$ws.Range("E30").Value = "  +0.01%  "

$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.55"
$ws.Range("E31").Value = "  -3.06%  "

$ws.Range("B32").Value = "VeChain"
$ws.Range("C32").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0455"
$ws.Range("E32").Value = "  +31.96%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.31"
$ws.Range("E33").Value = "  -0.45%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.80"
$ws.Range("E34").Value = "  +3.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0838"
$ws.Range("E35").Value = "  +1.43%  "

$ws.Range("E36").Value = "  -0.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.11"
$ws.Range("E37").Value = "  -0.44%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.79"
$ws.Range("E38").Value = "  -3.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.96"
$ws.Range("E39").Value = "  -0.91%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.23"
$ws.Range("E40").Value = "  -1.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.62"
$ws.Range("E41").Value = "  +10.70%  "

$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.47"
$ws.Range("E42").Value = "  -0.74%  "

$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.116"
$ws.Range("E43").Value = "  +1.87%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "126.01"
$ws.Range("E44").Value = "  -3.20%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.30"
$ws.Range("E45").Value = "  +0.88%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.33"
$ws.Range("E46").Value = "  -1.50%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.055.87"
$ws.Range("E47").Value = "  -2.84%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.30"
$ws.Range("E48").Value = "  +1.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.955"
$ws.Range("E49").Value = "  +9.35%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.58"
$ws.Range("E50").Value = "  -1.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.99"
$ws.Range("E51").Value = "  -0.15%  "
